# DataPool_v2.xlsx update:
# "Test 67 al 110 Factura Afecta y Test 01 al 35 Factura Exenta"
#
# - Row 164 (FA_0066): USUARIO (RUT) changes from 13712759-8 to 18215678-7
# - New rows 165-248 appended, each with columns A=TC, B=USUARIO, C=PASSWORD:
#     FA_0067 .. FA_0110  (Factura Afecta test cases)
#     FE_0001 .. FE_0040  (Factura Exenta test cases)
#   all sharing the same USUARIO/PASSWORD pair as the new RUT.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newUsuario = "18215678-7"
$password   = "Verity4.0"

# The existing FA_0066 row keeps its test-case name, only the RUT changes.
$ws.Range("B164").Value = $newUsuario

$row = 165

# Factura Afecta: FA_0067 .. FA_0110
for ($n = 67; $n -le 110; $n++) {
    $tc = "FA_{0:D4}" -f $n
    $ws.Cells.Item($row, 1).Value = $tc
    $ws.Cells.Item($row, 2).Value = $newUsuario
    $ws.Cells.Item($row, 3).Value = $password
    $row++
}

# Factura Exenta: FE_0001 .. FE_0040
for ($n = 1; $n -le 40; $n++) {
    $tc = "FE_{0:D4}" -f $n
    $ws.Cells.Item($row, 1).Value = $tc
    $ws.Cells.Item($row, 2).Value = $newUsuario
    $ws.Cells.Item($row, 3).Value = $password
    $row++
}

# Reflect the cursor position Excel would have left after typing the last rows.
$ws.Range("E246").Select()
